$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (sCs / Nppc / Npr3 / FAPs) with revised values ---
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.065932333333333
$ws.Range("H2").Value = 6.197797
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3190146666666667
$ws.Range("N2").Value = 0.957044
$ws.Range("O2").Value = 0.1466753431539481
$ws.Range("P2").Value = 0.1466753431539481
$ws.Range("Q2").Value = 0.6590627146742222
$ws.Range("R2").Value = 5.931564432068
$ws.Range("S2").Value = 0.1466753431539481
$ws.Range("T2").Value = 0.1466753431539481

# --- Add new row 3 (sCs / Nppc / Npr3 / ECs) ---
$ws.Range("A3").Value = "sCs"
$ws.Range("B3").Value = "Nppc"
$ws.Range("C3").Value = "Npr3"
$ws.Range("D3").Value = "ECs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.065932333333333
$ws.Range("H3").Value = 6.197797
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.855956666666666
$ws.Range("N3").Value = 5.567869999999999
$ws.Range("O3").Value = 0.8533246568460519
$ws.Range("P3").Value = 0.853324656846052
$ws.Range("Q3").Value = 3.834280886932221
$ws.Range("R3").Value = 34.50852798238999
$ws.Range("S3").Value = 0.8533246568460519
$ws.Range("T3").Value = 0.853324656846052
